# User Mobility workbook update:
#  - Action dropdown value changed from "delete" to "add" (row 2)
#  - New user row added: mhemaraju@auchan.com / add (row 3), with a
#    mailto hyperlink on the email cell matching the existing row's style
#  - Active selection left on B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change the Action from "delete" to "add"
$ws.Range("C2").Value = "add"

# Row 3: new user entry
$ws.Range("A3").Value = "mhemaraju@auchan.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:mhemaraju@auchan.com") | Out-Null
# Re-apply the same hyperlink formatting used by the existing entry (A2)
$ws.Range("A3").Style = $ws.Range("A2").Style

$ws.Range("C3").Value = "add"

# Leave the selection on B3, matching the saved workbook state
$ws.Range("B3").Select() | Out-Null
